$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new export timestamp (2024-07-18 09:32:57 -> 2024-07-19 10:29:17)
$ws.Name = "IClientBalance-20240719-102917-"

# All the "Data" column (G) values move forward one day: 45491 (2024-07-18) -> 45492 (2024-07-19)
$ws.Range("G2:G275").Value = 45492

# Row 58: Valor Devedor (D) drops to 0, Valor Credor (E) rises, Saldo (H) follows E
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 11178.24
$ws.Range("H58").Value = 11178.24

# Row 119: Valor Credor (E) and Saldo (H) updated
$ws.Range("E119").Value = 19617.75
$ws.Range("H119").Value = 19617.75

# Row 255: Valor Credor (E) and Saldo (H) updated
$ws.Range("E255").Value = 150606.84
$ws.Range("H255").Value = 150606.84

# Row 270: Valor Devedor (D) becomes negative, Saldo (H) recalculated (D+E)
$ws.Range("D270").Value = -110.23
$ws.Range("H270").Value = 523.66

# Update the active selection shown when the workbook is reopened
$ws.Range("K22").Select() | Out-Null
